# Edit corresponding to commit "Fri, Jun 26, 2020  5:04:44 AM":
#   1. Slide 16's table switches from the deck's custom "Table_0" style
#      ({155E508A-D569-486B-9768-0C488B1D81E0}) to the built-in table
#      style {C997E84D-80A6-4F13-AD7A-F9F5724BBE33}.
#   2. The presentation's theme colour palette changes from the
#      "Integral" palette to the stock "Office Theme" palette.

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 16 -------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
$table = $tableShape.Table
$table.ApplyStyle("{C997E84D-80A6-4F13-AD7A-F9F5724BBE33}")

# --- 2. Recolor the theme: Integral -> Office Theme -------------------
# Theme colour indices (MsoThemeColorSchemeIndex order):
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
# .RGB uses the usual COM BGR-packed integer, i.e. 0xBBGGRR for hex
# colour RRGGBB.
$colors = $p.Slides.Item(1).ThemeColorScheme
$colors.Item(1).RGB  = 0x000000   # dk1      000000
$colors.Item(2).RGB  = 0xFFFFFF   # lt1      FFFFFF
$colors.Item(3).RGB  = 0x6A5444   # dk2      44546A
$colors.Item(4).RGB  = 0xE6E6E7   # lt2      E7E6E6
$colors.Item(5).RGB  = 0xD59B5B   # accent1  5B9BD5
$colors.Item(6).RGB  = 0x317DED   # accent2  ED7D31
$colors.Item(7).RGB  = 0xA5A5A5   # accent3  A5A5A5
$colors.Item(8).RGB  = 0x00C0FF   # accent4  FFC000
$colors.Item(9).RGB  = 0xC47244   # accent5  4472C4
$colors.Item(10).RGB = 0x47AD70   # accent6  70AD47
$colors.Item(11).RGB = 0xC16305   # hlink    0563C1
$colors.Item(12).RGB = 0x724F95   # folHlink 954F72
